# Append 45 new device rows (102-146) to the bottom of the existing table.
# The new rows continue the existing 9-row (regcntr_id/machine_id) rotation
# pattern while device_id keeps incrementing sequentially from the last
# existing row (101 -> device_id 3000120, so the new rows start at 3000121).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 102
$startDevice = 3000121
$rowCount = 45

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    $cycle = $i % 9

    $ws.Cells.Item($row, 1).Value = 10002 + $cycle
    $ws.Cells.Item($row, 2).Value = 10021 + $cycle
    $ws.Cells.Item($row, 3).Value = $startDevice + $i
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin()"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
}

# The saved workbook's Page Setup was switched to portrait orientation.
$ws.PageSetup.Orientation = 1

# The saved selection sits on the empty rows right below the new data
# (i.e. the user selected the remainder of the sheet after entering data).
$null = $ws.Range("A147:XFD1048576").Select()
